$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "time_taken" in F1, copying the header style (bold,
# bordered, centered) from the neighboring E1 cell so the same style index
# is reused rather than creating a brand-new (duplicate) style entry.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate F2:F36 with the per-row "time_taken" timestamp metadata values.
# These are plain text values (no special style), mirroring the data cells.
$ws.Range("F2").Value = "2021-10-05 10:52:35.609147"
$ws.Range("F3").Value = "2021-10-05 10:52:35.609163"
$ws.Range("F4").Value = "2021-10-05 10:52:35.609167"
$ws.Range("F5").Value = "2021-10-05 10:52:35.609170"
$ws.Range("F6").Value = "2021-10-05 10:52:35.609173"
$ws.Range("F7").Value = "2021-10-05 10:52:35.609176"
$ws.Range("F8").Value = "2021-10-05 10:52:35.609178"
$ws.Range("F9").Value = "2021-10-05 10:52:35.609181"
$ws.Range("F10").Value = "2021-10-05 10:52:35.609184"
$ws.Range("F11").Value = "2021-10-05 10:52:35.609186"
$ws.Range("F12").Value = "2021-10-05 10:52:35.609190"
$ws.Range("F13").Value = "2021-10-05 10:52:35.609192"
$ws.Range("F14").Value = "2021-10-05 10:52:35.609195"
$ws.Range("F15").Value = "2021-10-05 10:52:35.609198"
$ws.Range("F16").Value = "2021-10-05 10:52:35.609201"
$ws.Range("F17").Value = "2021-10-05 10:52:35.609203"
$ws.Range("F18").Value = "2021-10-05 10:52:35.609206"
$ws.Range("F19").Value = "2021-10-05 10:52:35.609209"
$ws.Range("F20").Value = "2021-10-05 10:52:35.609211"
$ws.Range("F21").Value = "2021-10-05 10:52:35.609214"
$ws.Range("F22").Value = "2021-10-05 10:52:35.609216"
$ws.Range("F23").Value = "2021-10-05 10:52:35.609219"
$ws.Range("F24").Value = "2021-10-05 10:52:35.609222"
$ws.Range("F25").Value = "2021-10-05 10:52:35.609224"
$ws.Range("F26").Value = "2021-10-05 10:52:35.609227"
$ws.Range("F27").Value = "2021-10-05 10:52:35.609230"
$ws.Range("F28").Value = "2021-10-05 10:52:35.609233"
$ws.Range("F29").Value = "2021-10-05 10:52:35.609235"
$ws.Range("F30").Value = "2021-10-05 10:52:35.609238"
$ws.Range("F31").Value = "2021-10-05 10:52:35.609240"
$ws.Range("F32").Value = "2021-10-05 10:52:35.609243"
$ws.Range("F33").Value = "2021-10-05 10:52:35.609246"
$ws.Range("F34").Value = "2021-10-05 10:52:35.609249"
$ws.Range("F35").Value = "2021-10-05 10:52:35.609252"
$ws.Range("F36").Value = "2021-10-05 10:52:35.609254"
